$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Clear the "Title" value (previously "Medications")
$ws.Range("B5").Value = ""

# Update the "Date" value to the new publish timestamp
$ws.Range("B8").Value = "2024-06-04T08:55:54+00:00"
